# TC_71807 - Battery Standby / Alarm Load template update
# Move the "AlarmLoadingDetail/StandbyLoadingDetail" label row and the
# "Battery Alarm (A)/Battery Standby (A)" value row from columns S:T
# (row 7/8) up to F:G (row 1/2), and drop the now-unused S:T columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the two header/value pairs (with their formatting) into their new home.
$ws.Range("S7").Copy($ws.Range("F1"))
$ws.Range("T7").Copy($ws.Range("G1"))
$ws.Range("S8").Copy($ws.Range("F2"))
$ws.Range("T8").Copy($ws.Range("G2"))

# Remove the old cells now that the content has moved.
$ws.Range("S7:T8").Clear()

# Column G now carries the (longer) "StandbyLoadingDetail" style header,
# matching the width that column T used to have.
$ws.Columns.Item(7).ColumnWidth = 18.83

# Reset the view: no frozen/scrolled top-left cell, selection on the
# newly relocated block.
$win = $excel.ActiveWindow()
$win.ScrollRow = 1
$win.ScrollColumn = 1
$null = $ws.Range("F1:G2").Select()

Write-Host "Moved battery standby / alarm load template cells to F1:G2"
